$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the font styling previously applied to rows 2 and 3 (A2:D3) so they
# fall back to the default (unstyled) cell format, matching the rest of the
# data rows.
$ws.Range("A2:D3").Style = "Normal"

# Add the new "Spanish" language row.
$ws.Range("A8").Value = "spa"
$ws.Range("B8").Value = "Spanish"
$ws.Range("C8").Value = "Indo-European"
$ws.Range("D8").Value = "Spanish"

# Copy E7 (code/style for the "is_active" = TRUE column) into E8 so the new
# row keeps the same text-shared-string "TRUE" value and style as the rest
# of the table instead of being re-interpreted as a boolean.
$ws.Range("E7").Copy($ws.Range("E8"))

$ws.Range("A8").Select() | Out-Null
